# Handback report generation: update status text, populate "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns for the
# zh-cn and de-de localization sheets, and link the handback file names.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status moves from "Ready for handoff" to "Handed back: in sync with
#    en-US" everywhere it is shown (Overview status columns + per-language
#    Status column).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in the handback target file, handback xliff file name
#    and handback datetime for both rows.
# ---------------------------------------------------------------------------
$mdFileName   = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.md"
$mdHref       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1846cdb80f51c49ebc9488fc5fe1b9a7d45a4dcc/e2e/e03e8cad-5f2b-4759-a80d-0581a2aa856c.md"

$zhXlf = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.101179c73998b821a8504f720cbefac42762ec1d.zh-cn.xlf"
$deXlf = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.101179c73998b821a8504f720cbefac42762ec1d.de-de.xlf"

$zhHandbackTime = "2016-09-01 05:07:33"
$deHandbackTime = "2016-09-01 05:07:41"

$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Range("J2").Value = $zhXlf
$wsZhCn.Range("K2").Value = $zhHandbackTime

$wsZhCn.Range("I3").Value = $mdFileName
$wsZhCn.Range("J3").Value = $zhXlf
$wsZhCn.Range("K3").Value = $zhHandbackTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdHref, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdHref, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null

# ---------------------------------------------------------------------------
# 3. de-de sheet: same fields, using the de-de xliff + its own handback time.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Range("J2").Value = $deXlf
$wsDeDe.Range("K2").Value = $deHandbackTime

$wsDeDe.Range("I3").Value = $mdFileName
$wsDeDe.Range("J3").Value = $deXlf
$wsDeDe.Range("K3").Value = $deHandbackTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdHref, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdHref, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null

# ---------------------------------------------------------------------------
# 4. Widen the now-longer columns so the new text is fully visible (mirrors
#    Excel auto-fit behaviour after the wider handback status / file names
#    were written).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item("E").ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item("F").ColumnWidth = 29.9777047293527

$wsZhCn.Columns.Item("C").ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item("I").ColumnWidth = 40
$wsZhCn.Columns.Item("J").ColumnWidth = 40

$wsDeDe.Columns.Item("C").ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item("I").ColumnWidth = 40
$wsDeDe.Columns.Item("J").ColumnWidth = 40
